# Apply updated crypto price/volume data (and two row-position swaps)
# to match the refreshed "cryptos" listing, as produced by the daily
# GitHub Actions data-refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.711.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4935"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3001"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06775"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.918.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07328"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.212"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6752"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.678.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007974"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.164.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.384"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "198.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.331"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.674"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.961"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("E29").Value = "  +4.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.350"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.076"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05273"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7443"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("E35").Value = "  +1.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.714"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01854"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.719"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9300"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.082"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4500"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.77%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.78%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.62%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +23.99%  "

$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1400"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.00%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.060"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05901"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4039"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.89%  "
